$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InApkStringTable")

# Insert two new rows above row 274 to register the JellyFishGirl character
# (id / kor / eng columns), shifting every row below it down by 2.
$ws.Range("A274:D275").EntireRow.Insert()

# The inserted rows copy the formatting of row 273 (a wrapped description
# row), so row 274 (a single-line name row) needs its inherited wrap-text
# formatting cleared back to the default style.
$ws.Range("C274").WrapText = $false

# Column A (ids) first ...
$ws.Range("A274").Value = "CharName_JellyFishGirl"
$ws.Range("A275").Value = "CharDesc_JellyFishGirl"

# ... then column C (Korean text) ...
$ws.Range("C274").Value = "젤리피쉬걸"
$ws.Range("C275").Value = "젤리피쉬걸의 설명 우다다다`n곡사로 공격한다"
$ws.Range("C275").WrapText = $true
$ws.Rows(275).RowHeight = 49.5

# ... then column D (English text / translation placeholder).
$ws.Range("D274").Value = "JellyFIshGirl"
$ws.Range("D275").Formula = "=""In progress of translating…(""&ROW()&"")"""

# Duplicate-count helper formulas for the two new rows.
$ws.Range("B274").Formula = "=COUNTIF(A:A,A274)"
$ws.Range("B275").Formula = "=COUNTIF(A:A,A275)"

# The hidden _xlnm._FilterDatabase name tracked the InApkStringTable id column;
# extend it to cover the two freshly inserted rows.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = '=InApkStringTable!$B$1:$B$291'
